$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raluca")

$ws.Range("C26").Value = "https://creativeclass.carandache.com/products/lavaux-a-swiss-landscape"
$ws.Range("D26").Value = "19 CHF"
$ws.Range("C27").Value = "https://creativeclass.carandache.com/products/poetic-lily"
$ws.Range("B27").Value = "https://vhx.imgix.net/carandache/assets/0803cd32-cc66-42e9-98f6-246c5c876cb2.jpg?auto=format%2Ccompress&fit=crop&h=720&w=1280"
$ws.Range("A27").Value = "Curs Poetic Lily"
$ws.Range("B26").Value = "https://vhx.imgix.net/carandache/assets/6b347c23-113e-4831-af8f-fc01524383b1.jpg?auto=format%2Ccompress&fit=crop&h=720&w=1280"
$ws.Range("A26").Value = "Curs Lavaux"
$ws.Range("D27").Value = "19 CHF"

$ws.Range("A26").WrapText = $true
$ws.Range("A27").WrapText = $true

$ws.Range("D19").Select() | Out-Null
